$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to be written as text, even when the string looks
    # like a number (e.g. "530.68"), without leaving a residual explicit
    # cell style behind.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '69.552.42'
$ws.Range("E2").Value = '  +1.87%  '
$ws.Range("D3").Value = '3.932.65'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +0.21%  '
Set-TextValue "D5" '530.68'
$ws.Range("E5").Value = '  +8.61%  '
Set-TextValue "D6" '145.08'
$ws.Range("E6").Value = '  -1.25%  '
Set-TextValue "D7" '0.619'
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("E8").Value = '  -0.02%  '
Set-TextValue "D9" '0.729'
$ws.Range("E9").Value = '  +0.11%  '
Set-TextValue "D10" '0.173'
$ws.Range("E10").Value = '  +3.91%  '
Set-TextValue "D11" '0.0000340'
$ws.Range("E11").Value = '  -1.81%  '
Set-TextValue "D12" '42.68'
$ws.Range("E12").Value = '  -0.97%  '
Set-TextValue "D13" '10.41'
$ws.Range("E13").Value = '  -4.16%  '
$ws.Range("D14").Value = '4.565.56'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").Value = '3.938.19'
$ws.Range("E15").Value = '  +0.37%  '
Set-TextValue "D16" '14.03'
$ws.Range("E16").Value = '  -1.32%  '
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("E18").Value = '  +6.88%  '
Set-TextValue "D19" '19.84'
$ws.Range("E19").Value = '  -0.33%  '
$ws.Range("D20").Value = '69.399.85'
$ws.Range("E20").Value = '  +1.46%  '
Set-TextValue "D21" '431.87'
$ws.Range("E21").Value = '  -0.67%  '
$ws.Range("E22").Value = '  -4.13%  '
Set-TextValue "D23" '14.51'
$ws.Range("E23").Value = '  -2.48%  '
Set-TextValue "D24" '88.68'
$ws.Range("E24").Value = '  +1.00%  '
Set-TextValue "D25" '4.10'
$ws.Range("E25").Value = '  +13.42%  '
Set-TextValue "D26" '11.84'
$ws.Range("E26").Value = '  +3.79%  '
Set-TextValue "D27" '10.81'
$ws.Range("E27").Value = '  -4.14%  '
Set-TextValue "D28" '36.62'
$ws.Range("E28").Value = '  -4.02%  '
Set-TextValue "D29" '702.95'
$ws.Range("E29").Value = '  -3.11%  '
Set-TextValue "D30" '13.33'
$ws.Range("E30").Value = '  -3.38%  '
$ws.Range("E31").Value = '  -1.96%  '
$ws.Range("E32").Value = '  -1.90%  '
Set-TextValue "D33" '70.15'
$ws.Range("E33").Value = '  +16.19%  '
Set-TextValue "D34" '0.454'
$ws.Range("E34").Value = '  +12.29%  '
Set-TextValue "D35" '6.11'
$ws.Range("E35").Value = '  -2.89%  '
$ws.Range("D36").Value = '0.0₃0862'
$ws.Range("E36").Value = '  -1.64%  '
Set-TextValue "D37" '40.46'
$ws.Range("E37").Value = '  -2.99%  '
$ws.Range("E38").Value = '  +0.63%  '
Set-TextValue "D39" '0.998'
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("E40").Value = '  +0.01%  '
Set-TextValue "D41" '0.0482'
$ws.Range("E41").Value = '  +0.14%  '
Set-TextValue "D42" '2.82'
$ws.Range("E42").Value = '  -5.07%  '
Set-TextValue "D43" '3.11'
$ws.Range("E43").Value = '  +6.77%  '
$ws.Range("E44").Value = '  -5.13%  '
Set-TextValue "D45" '3.22'
$ws.Range("E45").Value = '  +14.28%  '
Set-TextValue "D46" '3.39'
$ws.Range("E46").Value = '  +2.46%  '
$ws.Range("E47").Value = '  +0.72%  '
$ws.Range("D48").Value = '0.0₆0357'
$ws.Range("E48").Value = '  +0.31%  '
Set-TextValue "D49" '3.33'
$ws.Range("E49").Value = '  -2.33%  '

# Rows 50/51: ARBITRUM and Monero swap positions, with updated prices/volumes
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D50" '144.93'
$ws.Range("E50").Value = '  +0.10%  '

$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D51" '2.09'
$ws.Range("E51").Value = '  -2.17%  '
